# Update "想去人数" (want-to-go count) figures across sheets to the
# newly generated values from gh-pages output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9529
$ws1.Range("F4").Value = 26
$ws1.Range("F5").Value = 530
$ws1.Range("F6").Value = 465

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 2

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9529
$ws4.Range("F4").Value = 26
$ws4.Range("F5").Value = 530
$ws4.Range("F6").Value = 2
$ws4.Range("F7").Value = 465
